$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C1").Value = 0.84260864398009372
$ws.Range("BO1").Value = 0.85636043113129801
$ws.Range("A2").Value = 0.72260391759216158
$ws.Range("T2").Value = 0.61427445211387766
$ws.Range("AV2").Value = 0.87223822229738079
$ws.Range("BP2").Value = 0.95561318637411496
$ws.Range("B3").Value = 0.72356658911305649
$ws.Range("D3").Value = 0.79999070243168657
$ws.Range("E3").Value = 0.71302708291378925
$ws.Range("F4").Value = 0.8983675389543273
$ws.Range("AB4").Value = 0.74203979960884769
$ws.Range("F5").Value = 0.85172336045321939
$ws.Range("AC5").Value = 0.95216126909715415
$ws.Range("AR5").Value = 0.70541513536557221
$ws.Range("AC6").Value = 0.60160439032789437
$ws.Range("I7").Value = 0.75648578985886827
$ws.Range("L7").Value = 0.97282442241229927
$ws.Range("BJ8").Value = 0.90051238204334483
$ws.Range("BM8").Value = 0.8855198391591157
$ws.Range("K9").Value = 0.84644733399443361
$ws.Range("J11").Value = 0.83129781019991444
$ws.Range("K12").Value = 0.84508296626863888
$ws.Range("M12").Value = 0.83821147928871764
$ws.Range("R12").Value = 0.59380287353695105
$ws.Range("A13").Value = 0.86185920866143073
$ws.Range("K14").Value = 0.66180478593018854
$ws.Range("O14").Value = 0.88896650025264612
$ws.Range("P15").Value = 0.98226062423690985
$ws.Range("BL15").Value = 0.8549205516105689
$ws.Range("N16").Value = 0.89692839569238503
$ws.Range("AA16").Value = 0.83162208819788452
$ws.Range("R17").Value = 0.99964222584123619
$ws.Range("S17").Value = 0.72672929385055873
$ws.Range("S18").Value = 0.82438113918328804
$ws.Range("U19").Value = 0.8665687270232536
$ws.Range("V19").Value = 0.95371502425060473
$ws.Range("R20").Value = 0.89001598649757696
$ws.Range("T21").Value = 0.99780883529992148
$ws.Range("V21").Value = 0.77693043260894568
$ws.Range("W22").Value = 0.68070208370088525
$ws.Range("U23").Value = 0.86015877671686036
$ws.Range("X23").Value = 0.72375191815885898
$ws.Range("Y23").Value = 0.96950247294022907
$ws.Range("V24").Value = 0.93975311116069782
$ws.Range("AA24").Value = 0.70745623425075799
$ws.Range("BC24").Value = 0.87076814040104034
$ws.Range("X25").Value = 0.91778065171644974
$ws.Range("Y27").Value = 0.99398638976933551
$ws.Range("Z27").Value = 0.87511676224257573
$ws.Range("AP27").Value = 0.8529050576631706
$ws.Range("J28").Value = 0.83053221842952274
$ws.Range("BO28").Value = 0.98299376254424675
$ws.Range("AA29").Value = 0.96759103683148129
$ws.Range("AB30").Value = 0.84829026348035597
$ws.Range("AE30").Value = 0.78572792093810817
$ws.Range("BJ30").Value = 0.79103429236839662
$ws.Range("AC31").Value = 0.94827288315378688
$ws.Range("AG31").Value = 0.84046613814725968
$ws.Range("N32").Value = 0.79603242229038229
$ws.Range("AE32").Value = 0.85136151652969816
$ws.Range("AD33").Value = 0.98825679689507684
$ws.Range("AI33").Value = 0.57917514608634235
$ws.Range("AN33").Value = 0.80041184017810574
$ws.Range("D34").Value = 0.79954751632413967
$ws.Range("Z34").Value = 0.97675724494863914
$ws.Range("AF34").Value = 0.63355736415809805
$ws.Range("BL34").Value = 0.66869500458060482
$ws.Range("AI36").Value = 0.93219071016708721
$ws.Range("AK36").Value = 0.93663981069834801
$ws.Range("AI37").Value = 0.93459166269701677
$ws.Range("AN37").Value = 0.77872357417198645
$ws.Range("AC38").Value = 0.87497678230139286
$ws.Range("AJ38").Value = 0.89341876792731156
$ws.Range("AK38").Value = 0.85916509736882918
$ws.Range("AL39").Value = 0.69411331283705313
$ws.Range("AO39").Value = 0.97258421004058304
$ws.Range("AM40").Value = 0.67766360062747544
$ws.Range("AO40").Value = 0.97605478961540793
$ws.Range("AP41").Value = 0.72794956439635605
$ws.Range("AD42").Value = 0.94924794463396056
$ws.Range("Y44").Value = 0.94272717175784759
$ws.Range("AQ44").Value = 0.71145965476602191
$ws.Range("BH44").Value = 0.58049612173433907
$ws.Range("AQ45").Value = 0.92885848537073867
$ws.Range("AR45").Value = 0.87392239094461055
$ws.Range("AT45").Value = 0.73114890095264629
$ws.Range("AU46").Value = 0.65289015634703751
$ws.Range("AS47").Value = 0.98520965819974937
$ws.Range("BJ47").Value = 0.94088040055697775
$ws.Range("AT48").Value = 0.55065647201589996
$ws.Range("AU49").Value = 0.95362748929474039
$ws.Range("AX49").Value = 0.95155103953600251
$ws.Range("AV50").Value = 0.88097306901445971
$ws.Range("AW51").Value = 0.61845144434238941
$ws.Range("BO51").Value = 0.7554658142677263
$ws.Range("L52").Value = 0.76965640273485958
$ws.Range("AX52").Value = 0.58129983349179215
$ws.Range("AY52").Value = 0.56703438554098073
$ws.Range("BG52").Value = 0.96812176343446366
$ws.Range("F54").Value = 0.94749751427930873
$ws.Range("BA54").Value = 0.75193602671243065
$ws.Range("BD54").Value = 0.57962842363541656
$ws.Range("BE54").Value = 0.99142507783417821
$ws.Range("BE55").Value = 0.75915016096332066
$ws.Range("H56").Value = 0.88567909489156571
$ws.Range("BC56").Value = 0.85955385617861935
$ws.Range("BF56").Value = 0.94109512734810719
$ws.Range("I57").Value = 0.82028882159353289
$ws.Range("M57").Value = 0.92336186064949777
$ws.Range("BF57").Value = 0.72703248535979936
$ws.Range("BH58").Value = 0.7346004793022074
$ws.Range("E59").Value = 0.71493470539019421
$ws.Range("BE59").Value = 0.81550293376158312
$ws.Range("BF59").Value = 0.78838978964842688
$ws.Range("BI59").Value = 0.98562710734646819
$ws.Range("AO60").Value = 0.71610996364692947
$ws.Range("BA60").Value = 0.90788557587476837
$ws.Range("BG60").Value = 0.76360251059088502
$ws.Range("BI60").Value = 0.90149946610560516
$ws.Range("BI62").Value = 0.88551929673017094
$ws.Range("BK62").Value = 0.82368897286427845
$ws.Range("AV63").Value = 0.94130066348426977
$ws.Range("BI63").Value = 0.99208313057833286
$ws.Range("BM63").Value = 0.97314666631242952
$ws.Range("BJ64").Value = 0.7421656822326218
$ws.Range("BK64").Value = 0.993576930826088
$ws.Range("P65").Value = 0.98032227918675374
$ws.Range("BL65").Value = 0.72962046516951529
$ws.Range("AX66").Value = 0.76978201261204937
$ws.Range("BL66").Value = 0.72010121463328525
$ws.Range("BM66").Value = 0.76392402036385354
$ws.Range("BO66").Value = 0.95782195192726438
$ws.Range("BP66").Value = 0.86181332067043548
$ws.Range("BM67").Value = 0.76804671370495914
$ws.Range("BP67").Value = 0.92851256865122123
$ws.Range("AV68").Value = 0.91011386107838921
